# Commit: "remove the .1 from BX248355.1 in the BX248355.fasta headers and metadata"
#
# The sheet originally listed the BX248355 assembly as four synthetic
# "segments" (rows 4-6: BX248355.1-segment2/3/4). The edit collapses that
# back down to a single BX248355 row (row 4) whose identifiers drop the
# ".1" suffix, and removes the now-redundant segment rows (5 and 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diphtheria_Metadata")

# Drop the extra segment rows (work bottom-up so row numbers above stay put).
$ws.Rows(6).Delete()
$ws.Rows(5).Delete()

# Rewrite the remaining BX248355 row's identifiers without the ".1"/segment suffix.
$ws.Range("B4").Value = "BX248355.1"
$ws.Range("K4").Value = "BX248355.1"
$ws.Range("AF4").Value = "BX248355.1"
$ws.Range("A4").Value = "BX248355"

# Leave the cursor where the saved file's view last had it.
$ws.Range("AN4").Select() | Out-Null
